$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 1 (headers) ---
# Shift existing headers right: B1,C1,D1,E1,F1 -> D1,E1,F1,G1,H1
$ws.Range("D1").Value = "Name"
$ws.Range("E1").Value = "Static_Path"
$ws.Range("F1").Value = "VMM_Domain"
$ws.Range("G1").Value = "BD"
$ws.Range("H1").Value = "Contract"
# New header cells
$ws.Range("B1").Value = "Tenant"
$ws.Range("C1").Value = "ANP"

# --- Row 2 (existing Create_EPG entry, now shifted) ---
$ws.Range("A2").Value = "Create_EPG"
$ws.Range("B2").Value = "Prod"
$ws.Range("C2").Value = "e7_mgmt"
$ws.Range("D2").Value = "Storage_Mgmt"
$ws.Range("E2").Value = $null
$ws.Range("F2").Value = "e7vmw1vic02"
$ws.Range("G2").Value = "Storage_Mgmt_BD"
$ws.Range("H2").Value = "e7_L3_Out"

# --- Row 3 (new Create_EPG entry) ---
$ws.Range("A3").Value = "Create_EPG"
$ws.Range("B3").Value = "Prod"
$ws.Range("C3").Value = "e7_mgmt"
$ws.Range("D3").Value = "Network_Mgmt"
$ws.Range("F3").Value = "e7vmw1vic02"
$ws.Range("G3").Value = "Network_Mgmt"
$ws.Range("H3").Value = "e7_L3_Out"

# --- Row 5 (Create_BD headers) ---
$ws.Range("B5").Value = "Tenant"
$ws.Range("C5").Value = "Name"
$ws.Range("D5").Value = "VRF"
$ws.Range("E5").Value = "Subnet"
$ws.Range("F5").Value = "advertise (yes/no)"
$ws.Range("G5").Value = "L3 out"

# --- Row 6 (Create_BD entry, replaces old row 4 A4=Create_BD) ---
$ws.Range("A4").Value = $null
$ws.Range("A6").Value = "Create_BD"
$ws.Range("B6").Value = "Prod"
$ws.Range("C6").Value = "Storage_Mgmt_BD"
$ws.Range("D6").Value = "Prod"
$ws.Range("E6").Value = "10.207.50.1/24"
$ws.Range("F6").Value = "yes"

# --- Column widths ---
# Column A already has the correct width (21.6640625) from the original file, leave untouched.
# B and C take on the same width as column A.
$ws.Range("B:C").ColumnWidth = 20.833333333333332
$ws.Range("D:D").ColumnWidth = 17.666666666666668
$ws.Range("E:E").ColumnWidth = 19.666666666666668
$ws.Range("F:F").ColumnWidth = 17.666666666666668
$ws.Range("G:G").ColumnWidth = 25.5

# --- Selection ---
$ws.Range("G6").Select()

# --- Window height ---
$excel.ActiveWindow.Height = 16240
